$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings must be introduced in the same order the final workbook
# expects them to appear in the shared-string table.
# Row 14 (header row): add "No of Blocks" header in column K
$ws.Range("K14").Value = "No of Blocks"

# Row 16: block-count value "3" (new string), then the renamed test case/class
$ws.Range("K16").Value = "'3"
$ws.Range("D15").Value = "verifyCreateRoster"
$ws.Range("D16").Value = "updateCreateRoster"
$ws.Range("C15").Value = "rosterSetupTest"
$ws.Range("C16").Value = "rosterSetupTest"

# Row 15: block-count value "2" (reuses existing shared string)
$ws.Range("K15").Value = "'2"

# Row 16: clear old 14-day pattern columns
$ws.Range("H16").ClearContents()
$ws.Range("I16").ClearContents()
$ws.Range("J16").ClearContents()

# Update the view: scroll so column B is the left-most visible column, and move the selection
$ws.Range("H21").Select()
$excel.ActiveWindow.ScrollColumn = 2
